$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextValue "D2" "62.414.66"
Set-TextValue "E2" "  +1.17%  "
Set-TextValue "D3" "3.000.67"
Set-TextValue "E3" "  -0.44%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "544.76"
Set-TextValue "E5" "  -1.02%  "
Set-TextValue "D6" "139.21"
Set-TextValue "E6" "  +3.38%  "
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "D8" "2.995.85"
Set-TextValue "E8" "  -0.43%  "
Set-TextValue "D9" "0.489"
Set-TextValue "E9" "  -1.58%  "
Set-TextValue "D10" "6.77"
Set-TextValue "E10" "  +11.80%  "
Set-TextValue "D11" "0.148"
Set-TextValue "E11" "  -0.07%  "
Set-TextValue "D12" "0.447"
Set-TextValue "E12" "  -0.59%  "
Set-TextValue "E13" "  -0.50%  "
Set-TextValue "D14" "34.02"
Set-TextValue "E14" "  -0.80%  "
Set-TextValue "D15" "3.481.20"
Set-TextValue "E15" "  -0.84%  "
Set-TextValue "D16" "62.443.10"
Set-TextValue "E16" "  +1.01%  "
Set-TextValue "D17" "3.004.10"
Set-TextValue "E17" "  -0.58%  "
Set-TextValue "E18" "  -1.91%  "
Set-TextValue "D19" "6.57"
Set-TextValue "E19" "  -1.38%  "
Set-TextValue "D20" "469.65"
Set-TextValue "E20" "  -0.57%  "
Set-TextValue "D21" "13.41"
Set-TextValue "E21" "  +1.16%  "
Set-TextValue "D22" "0.653"
Set-TextValue "E22" "  -3.15%  "
Set-TextValue "D23" "7.19"
Set-TextValue "E23" "  +1.57%  "
Set-TextValue "E24" "  -0.80%  "
Set-TextValue "D25" "12.61"
Set-TextValue "E25" "  +4.09%  "
Set-TextValue "E26" "  -0.08%  "
Set-TextValue "E27" "  +0.10%  "
Set-TextValue "D28" "7.62"
Set-TextValue "E28" "  -2.43%  "
Set-TextValue "E29" "  +5.52%  "
Set-TextValue "E30" "  +0.06%  "
Set-TextValue "D31" "25.43"
Set-TextValue "E31" "  -1.09%  "
Set-TextValue "E32" "  -2.22%  "
Set-TextValue "E33" "  +2.20%  "
Set-TextValue "D34" "5.56"
Set-TextValue "E34" "  +1.72%  "
Set-TextValue "D35" "54.63"
Set-TextValue "E35" "  -1.36%  "
Set-TextValue "D36" "5.83"
Set-TextValue "E36" "  -1.32%  "
Set-TextValue "D37" "450.08"
Set-TextValue "E37" "  -1.74%  "
Set-TextValue "D38" "0.0810"
Set-TextValue "E38" "  +1.71%  "
Set-TextValue "D39" "0.0391"
Set-TextValue "E39" "  +2.00%  "
Set-TextValue "D40" "2.958.52"
Set-TextValue "E40" "  -8.00%  "
Set-TextValue "D42" "8.06"
Set-TextValue "E42" "  -1.18%  "
Set-TextValue "E43" "  +5.28%  "
Set-TextValue "D44" "26.79"
Set-TextValue "E44" "  +2.78%  "
Set-TextValue "D46" "0.247"
Set-TextValue "E46" "  +1.13%  "
Set-TextValue "E47" "  +0.70%  "
Set-TextValue "E48" "  +0.69%  "
Set-TextValue "D49" "115.29"
Set-TextValue "E49" "  -2.21%  "
Set-TextValue "D50" "0.0₃0496"
Set-TextValue "E50" "  +0.41%  "
Set-TextValue "D51" "2.01"
Set-TextValue "E51" "  -0.45%  "
